# "Guardar" action: append the current form values as a new row at the
# bottom of the parts list, right under the last used row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Form field values (blank text-box fields are saved as empty strings).
$refComponente = "11113"
$frontal = "1"
$lateralDer = ""
$lateralIzq = "3"
$powerReset = ""
$ledsFrontales = ""
$varios = "flaikers"
$protecciones = ""

$formRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 8))
# Keep every field as plain text (matches how the rest of the sheet stores
# its values) instead of letting Excel auto-convert numeric-looking text.
$formRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = $refComponente
$ws.Cells.Item($newRow, 2).Value = $frontal
$ws.Cells.Item($newRow, 3).Value = $lateralDer
$ws.Cells.Item($newRow, 4).Value = $lateralIzq
$ws.Cells.Item($newRow, 5).Value = $powerReset
$ws.Cells.Item($newRow, 6).Value = $ledsFrontales
$ws.Cells.Item($newRow, 7).Value = $varios
$ws.Cells.Item($newRow, 8).Value = $protecciones
